$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate rows 2-12 with the full card table (row 2 is overwritten with new
# card data; rows 3-12 are newly added cards).

# Row 2: Osgiv; the Rconstructor
$ws.Range("A2").Value = "Osgiv; the Rconstructor"
$ws.Range("B2").Value = "4 any; 1 Red and White"
$ws.Range("C2").Value = "Legendary Creature - Giant Artificer"
$ws.Range("D2").Value = "Vigilance; Cost 1 and Sacrifice and artifiact: Target creature yo ucontrol gets +2/+0 until the end of turn. Cost x and tap; Excile an artifact card with mana value X from your graveyard: Create two tokens that are copies of the exiled card. Activate only as a sorcery "
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = 4

# Row 3: Dispeller's Capsule
$ws.Range("A3").Value = "Dispeller's Capsule"
$ws.Range("B3").Value = "1 White"
$ws.Range("C3").Value = "Artifact"
$ws.Range("D3").Value = "2 any and 1 white and tap; Sacrifice Dispeller's Capsule: Destroy target artifact or enchantment."
$ws.Range("E3").Value = "n/a"
$ws.Range("F3").Value = "n/a"

# Row 4: Mycosynth Wellspring
$ws.Range("A4").Value = "Mycosynth Wellspring"
$ws.Range("B4").Value = "2 any"
$ws.Range("C4").Value = "Artifact"
$ws.Range("D4").Value = "When Mycosynth Wellspring enters the battlefield or is put inot a graveyard from the battlefield; you may search your library for a basic land card; reveal it; put it into your hand; then shuffle."
$ws.Range("E4").Value = "n/a"
$ws.Range("F4").Value = "n/a"

# Row 5: Bronze Guardian
$ws.Range("A5").Value = "Bronze Guardian"
$ws.Range("B5").Value = "4 any; 1 White"
$ws.Range("C5").Value = "Artifact Creature - Golem"
$ws.Range("D5").Value = "Double strike; Ward 2; Other artifacts you control have ward 2. Bronze Guardian's power is equal to the number of artifacts you control."
$ws.Range("E5").Value = "*"
$ws.Range("F5").Value = 5

# Row 6: Steel Hellkite
$ws.Range("A6").Value = "Steel Hellkite"
$ws.Range("B6").Value = "6 any"
$ws.Range("C6").Value = "Artifact Creature - Dragon"
$ws.Range("D6").Value = "Flying; 2 any: Steel Hellkite gets +1/+0 until the end of turn. X any: Destroy each nonland permanent with mana value X whose controller was dealt comabt damage by Steel Hellkite this turn. Activate only once each turn."
$ws.Range("E6").Value = 5
$ws.Range("F6").Value = 5

# Row 7: Angle of the Ruins
$ws.Range("A7").Value = "Angle of the Ruins"
$ws.Range("B7").Value = "5 any; 2 White"
$ws.Range("C7").Value = "Artifact Creatrue - Angle"
$ws.Range("D7").Value = "Flying; When Angel of the Ruins enters the battlefield; exile up to two target artifacts and/or enchantments. Plainscycleing: 2 Any (Discard this card: Search your library for a plains card; reveal it;put it into your hand; then shuffle.)"
$ws.Range("E7").Value = 5
$ws.Range("F7").Value = 7

# Row 8: Cursed Mirror
$ws.Range("A8").Value = "Cursed Mirror"
$ws.Range("B8").Value = "2 Any; 1 Red"
$ws.Range("C8").Value = "Artifact"
$ws.Range("D8").Value = "tap: Add 1 Red. As Cursed Mirror enters the battlefield; you may have it become a copy of any creature on the battlefield until end of turn; except it has haste."
$ws.Range("E8").Value = "n/a"
$ws.Range("F8").Value = "n/a"

# Row 9: Monologue Tax
$ws.Range("A9").Value = "Monologue Tax"
$ws.Range("B9").Value = "2 Any; 1 White"
$ws.Range("C9").Value = "Enchantment"
$ws.Range("D9").Value = "Whenever an opponent casts their second spell each turn; you create a Treasure token."
$ws.Range("E9").Value = "n/a"
$ws.Range("F9").Value = "n/a"

# Row 10: Daretti; Scrap Savant
$ws.Range("A10").Value = "Daretti; Scrap Savant"
$ws.Range("B10").Value = "3 Any; 1 Red"
$ws.Range("C10").Value = "Ledendary Planeswalker - Daretti"
$ws.Range("D10").Value = "|+2 toughness: Discart up to two cards; then draw that many cards. -2: Sacrifice an artifact. If you do; return target artifact card from your graveyard to the battlefeild. -10: You get an emblem with ""Whenever an artifact is put into your graveyard from the battlefield; return that card to the battlefield at the beginning of the next end step."" Daretti; Scrap Savant can be your commander.C13"
$ws.Range("E10").Value = "n/a"
$ws.Range("F10").Value = 3

# Row 11: Boros Locket
$ws.Range("A11").Value = "Boros Locket"
$ws.Range("B11").Value = "3 Any; 1 Red"
$ws.Range("C11").Value = "Artifact"
$ws.Range("D11").Value = "tap: Add 1 Red and 1 White. 4 Red or White; tap Scarifice Boros Locket: Draw two cards."
$ws.Range("E11").Value = "n/a"
$ws.Range("F11").Value = "n/a"

# Row 12: Secret Rendezvous
$ws.Range("A12").Value = "Secret Rendezvous"
$ws.Range("B12").Value = "1 Any; 2 White"
$ws.Range("C12").Value = "Sorcery"
$ws.Range("D12").Value = "You and target opponent each draw three cards"
$ws.Range("E12").Value = "n/a"
$ws.Range("F12").Value = "n/a"

# Update the active cell selection to match the saved workbook state (A2).
$ws.Range("A2").Select() | Out-Null

